$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: expected_delivery becomes a real date value (was text "2014-07-31")
$ws.Range("E2").Value = 41851
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"
$ws.Range("F2").Value = 32345
$ws.Range("H2").Value = "test-sku-6"

# Row 3
$ws.Range("H3").Value = "test-sku-7"

# Row 4: expected_delivery becomes a real date value (was text "2014-08-12")
$ws.Range("E4").Value = 41863
$ws.Range("E4").NumberFormat = "yyyy-mm-dd"
$ws.Range("F4").Value = 32346
$ws.Range("H4").Value = "test-sku-8"

# Row 5
$ws.Range("H5").Value = "test-sku-9"

# Update selection to match author's last cursor position
$ws.Range("G14").Select()
